# PARTICIPANTS_TEMPLATE.xlsx edit
# Renames the generic header row to "Participant ..." / "Primary Delegate ..."
# phrasing, refreshes several helper-text cells in row 2, turns on WrapText for
# the whole header block (rows 1-2), fixes column A / G widths and row 1
# height, restyles Q2 to match the real "Hyperlink" look of G2, and updates
# the window zoom / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 1 (header labels) - plain text swaps only; font/fill stay as-is.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Participant First Name"
$ws.Range("B1").Value = "Participant Middle Name"
$ws.Range("C1").Value = "Participant Last Name"
$ws.Range("D1").Value = "Participant Nickname"
$ws.Range("E1").Value = "Participant Sex"
$ws.Range("F1").Value = "Participant Date of Birth"
$ws.Range("G1").Value = "Participant Email Address"
$ws.Range("H1").Value = "Participant Phone"
$ws.Range("J1").Value = "Participant Alternative Phone"
$ws.Range("R1").Value = "Primary Delegate Phone Number"
$ws.Range("S1").Value = "Primary Delegate Phone Type"

# ---------------------------------------------------------------------------
# 2. Row 2 (helper / instructional text)
# ---------------------------------------------------------------------------
$ws.Range("G2").Value = "youremail@yourdomain.com If Participant is a minor the email address will not be stored"
$ws.Range("R2").Value = "Format per country 919-555-1212."
$ws.Range("Q2").Value = "youremail@yourdomain.com   Mandatory for a minor participant"

# H2: two runs - plain lead-in, bold instruction.
$h2Text = "Format per country 919-555-1212       For a minor participant enter the participant delegate's phone number"
$ws.Range("H2").Value = $h2Text
$h2BoldStart = 39
$h2BoldLen = $h2Text.Length - $h2BoldStart + 1
$ws.Range("H2").Characters($h2BoldStart, $h2BoldLen).Font.Bold = $true

# O2: "Text.         " (plain) + "Mandatory for a minor participant" (bold)
$o2Text = "Text.         Mandatory for a minor participant"
$ws.Range("O2").Value = $o2Text
$o2BoldStart = 15
$o2BoldLen = $o2Text.Length - $o2BoldStart + 1
$ws.Range("O2").Characters($o2BoldStart, $o2BoldLen).Font.Bold = $true

# P2: "Text.          " (plain, one more trailing space than O2) + bold tail
$p2Text = "Text.          Mandatory for a minor participant"
$ws.Range("P2").Value = $p2Text
$p2BoldStart = 16
$p2BoldLen = $p2Text.Length - $p2BoldStart + 1
$ws.Range("P2").Characters($p2BoldStart, $p2BoldLen).Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. Wrap text across the whole populated header block.
# ---------------------------------------------------------------------------
$ws.Range("A1:S2").WrapText = $true

# ---------------------------------------------------------------------------
# 4. Q2 was plain bold black; it is really a mailto: hyperlink like G2, so
#    restyle it to match the regular Excel "Hyperlink" look (blue, underlined,
#    not bold).
# ---------------------------------------------------------------------------
$q2Font = $ws.Range("Q2").Font
$q2Font.Bold = $false
$q2Font.Underline = 2
$q2Font.Color = $ws.Range("G2").Font.Color

# ---------------------------------------------------------------------------
# 5. Column width / row height tweaks.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.25   # -> width 13.08 ish (was 11)
$ws.Columns.Item(7).ColumnWidth = 20.33   # -> width 21.16 ish (was 20.66)
$ws.Rows.Item(1).RowHeight = 46.5

# ---------------------------------------------------------------------------
# 6. Window: zoom to 100%, scroll so column K is leftmost, select Q1.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 100
$ws.Range("Q1").Select()
$win.ScrollColumn = 11
$win.ScrollRow = 1
